$wb = $excel.ActiveWorkbook

# --- Notes sheet: update the "specific issue" note text -----------------
$wsNotes = $wb.Worksheets.Item("Notes")
$wsNotes.Range("A3").Value = "Specific issue: variant_num exceeds total_num"

# --- studies sheet: lowercase the study_ID header to study_id -----------
$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Range("A1").Value = "study_id"
$wsStudies.Range("A2").Select()

# --- surveys sheet: lowercase/relabel headers + explicit black header font --
$wsSurveys = $wb.Worksheets.Item("surveys")
$wsSurveys.Range("B1").Value = "survey_id"
$wsSurveys.Range("E1").Value = "latitude"
$wsSurveys.Range("F1").Value = "longitude"
$wsSurveys.Range("A1:K1").Font.Color = 0
$wsSurveys.Range("A1:K1").Select()

# --- counts sheet: remove duplicated row, flag variant_num > total_num --
$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Rows.Item(3).Delete()
$wsCounts.Range("C2").Value = 11
$wsCounts.Range("C3").Select()
